$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "43.406.19"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.274.67"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "123.28"
$ws.Range("E5").Value = "  +5.67%  "
$ws.Range("D6").Value = "266.42"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  +2.23%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "48.09"
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "9.23"
$ws.Range("E12").Value = "  +3.48%  "
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "2.617.71"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "2.268.65"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "43.529.91"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "72.26"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "2.43"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "235.40"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -1.76%  "
$ws.Range("E25").Value = "  -4.49%  "
$ws.Range("D26").Value = "12.00"
$ws.Range("E26").Value = "  +2.59%  "
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("D28").Value = "42.62"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "172.82"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  -2.16%  "
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("E36").Value = "  +12.45%  "
$ws.Range("D37").Value = "0.0375"
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "2.56"
$ws.Range("E40").Value = "  +5.14%  "
$ws.Range("E41").Value = "  -3.94%  "
$ws.Range("D42").Value = "73.79"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "5.68"
$ws.Range("E46").Value = "  -11.41%  "
$ws.Range("D47").Value = "73.98"
$ws.Range("E47").Value = "  +36.95%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "101.59"
$ws.Range("E51").Value = "  -1.41%  "
